# add XS breakpoints to metadata
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("github")

# --- update the existing breakpoints value for the XS rows (G16/G17) ---
# was "9055:10447" -> now "9056:10448"
$ws.Range("G16").Value = "9056:10448"
$ws.Range("G17").Value = "9056:10448"

# --- populate the clades_regions column (F16/F17) with the new rich text ---
$richText = "210:10029|Delta/21J,10449:29742|Omicron/BA.1/21K"
$ws.Range("F16").Value = $richText
$ws.Range("F17").Value = $richText

foreach ($addr in @("F16", "F17")) {
    $cell = $ws.Range($addr)

    # "Delta/21J" -> bold red
    $deltaRun = $cell.Characters(11, 9)
    $deltaRun.Font.Name = "Calibri"
    $deltaRun.Font.Size = 11
    $deltaRun.Font.Bold = $true
    $deltaRun.Font.Color = 255

    # ",10449:29742|" -> normal weight, default color
    $midRun = $cell.Characters(20, 13)
    $midRun.Font.Name = "Calibri"
    $midRun.Font.Size = 11
    $midRun.Font.Bold = $false
    $midRun.Font.Color = 0

    # "Omicron/BA.1/21K" -> bold green
    $omicronRun = $cell.Characters(33, 16)
    $omicronRun.Font.Name = "Calibri"
    $omicronRun.Font.Size = 11
    $omicronRun.Font.Bold = $true
    $omicronRun.Font.Color = 5287936
}

# --- sheet view: scroll pane + selected cell moved ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F14").Select()

Write-Output "applied XS breakpoints edit"
